$wb = $excel.ActiveWorkbook

# Update the Fibonacci retracement base value (Low), which cascades through
# all of the retracement/extension formulas on this sheet as well as the
# linked cells on the Povit sheet (I34:I36 reference Fibonacci!E18/E17/E16).
$wsFib = $wb.Worksheets.Item("Fibonacci")
$wsFib.Range("E9").Value = 10304.6
$wsFib.Range("E16").Select() | Out-Null

# Restore the Povit sheet as the active/selected sheet, matching the cursor
# position left behind after reviewing the recalculated Fibonacci numbers.
$wsPivot = $wb.Worksheets.Item("Povit")
$wsPivot.Range("I36").Select() | Out-Null
